$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("produk")

# Row height changes on the produk sheet
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(16).RowHeight = 30

# Copy A1:I20 from produk into a new worksheet placed right after "produk"
$src = $ws.Range("A1:I20")
$src.Copy()
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$newSheet.Range("A1").PasteSpecial(-4163)
$newSheet.Range("A1").PasteSpecial(-4122)

$newSheet.Range("H4").Select()
